# Updated cryptos list on Thu Nov 30 02:17:38 UTC 2023 with GitHub Actions
#
# Refresh the crypto prices / 1h volume percentages (and, for a few rows whose
# ranking order swapped, the coin name + link) in the "cryptos" worksheet.
#
# All of the Price (D) and Volume(1h) (E) columns are stored as *text* in the
# workbook (values like "37.809.05" aren't valid numbers, and the percentages
# keep their surrounding padding spaces), so we force the target range to
# Text format before writing, then restore the cell style to "Normal" so we
# don't leave a stray number format attached to the cells (matches how the
# original file has these cells with no explicit style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:E51")
$rng.NumberFormat = "@"

# Bitcoin
$ws.Range("D2").Value = '37.744.43'
$ws.Range("E2").Value = '  -0.64%  '

# Ethereum
$ws.Range("D3").Value = '2.027.67'
$ws.Range("E3").Value = '  -1.14%  '

# BNB
$ws.Range("D5").Value = '226.67'
$ws.Range("E5").Value = '  -1.29%  '

# XRP
$ws.Range("E6").Value = '  +0.11%  '

# Solana
$ws.Range("D7").Value = '59.79'
$ws.Range("E7").Value = '  +1.87%  '

# USDC
$ws.Range("E8").Value = '  +0.02%  '

# Cardano
$ws.Range("E9").Value = '  -1.49%  '

# Dogecoin
$ws.Range("D10").Value = '0.0816'
$ws.Range("E10").Value = '  +0.53%  '

# TRON
$ws.Range("E11").Value = '  +0.11%  '

# WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '2.326.52'
$ws.Range("E12").Value = '  -1.15%  '

# Chainlink
$ws.Range("D13").Value = '14.44'
$ws.Range("E13").Value = '  -1.02%  '

# Avalanche
$ws.Range("D14").Value = '20.96'
$ws.Range("E14").Value = '  +0.77%  '

# Polygon
$ws.Range("D15").Value = '0.757'
$ws.Range("E15").Value = '  +0.41%  '

# Polkadot
$ws.Range("D16").Value = '5.15'
$ws.Range("E16").Value = '  -2.84%  '

# WrappedEther
$ws.Range("D17").Value = '2.021.61'
$ws.Range("E17").Value = '  -0.96%  '

# WrappedBTC
$ws.Range("D18").Value = '37.681.53'
$ws.Range("E18").Value = '  -0.54%  '

# Row 19 was Uniswap, is now Litecoin (ranking order swapped with row 20)
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").Value = '69.73'
$ws.Range("E19").Value = '  -0.20%  '

# Row 20 was Litecoin, is now Uniswap
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '5.94'
$ws.Range("E20").Value = '  -5.37%  '

# ShibaInu
$ws.Range("D21").Value = '0.0₃0823'
$ws.Range("E21").Value = '  -1.61%  '

# BitcoinCash
$ws.Range("D22").Value = '223.98'
$ws.Range("E22").Value = '  -0.41%  '

# Dai
$ws.Range("E23").Value = '  -0.01%  '

# Toncoin
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  -1.14%  '

# PancakeSwap
$ws.Range("E25").Value = '  +0.20%  '

# Row 26 was Cosmos, is now Monero (ranking order swapped with row 27)
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '165.75'
$ws.Range("E26").Value = '  -0.43%  '

# Row 27 was Monero, is now Cosmos
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '9.25'
$ws.Range("E27").Value = '  -0.74%  '

# Kaspa
$ws.Range("E28").Value = '  -4.17%  '

# EthereumClassic
$ws.Range("D29").Value = '18.81'
$ws.Range("E29").Value = '  -1.23%  '

# ImmutableX
$ws.Range("D30").Value = '1.26'
$ws.Range("E30").Value = '  -5.26%  '

# Stellar
$ws.Range("E31").Value = '  +0.77%  '

# WEMIXToken
$ws.Range("D32").Value = '2.17'
$ws.Range("E32").Value = '  +6.30%  '

# Filecoin
$ws.Range("D33").Value = '4.38'
$ws.Range("E33").Value = '  -3.20%  '

# Row 34 was InternetComputer(DFINITY), is now Hedera (ranking order swapped with row 35)
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0603'
$ws.Range("E34").Value = '  -1.58%  '

# Row 35 was Hedera, is now InternetComputer(DFINITY)
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '4.48'
$ws.Range("E35").Value = '  -2.40%  '

# THORChain
$ws.Range("E36").Value = '  +4.62%  '

# LidoDAOToken
$ws.Range("D37").Value = '2.25'
$ws.Range("E37").Value = '  -2.90%  '

# RenderToken
$ws.Range("D38").Value = '3.22'
$ws.Range("E38").Value = '  -2.28%  '

# BinanceUSD
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  -0.05%  '

# Maker
$ws.Range("D40").Value = '1.519.53'
$ws.Range("E40").Value = '  +2.52%  '

# VeChain
$ws.Range("E41").Value = '  -0.37%  '

# InjectiveProtocol
$ws.Range("D42").Value = '16.79'

# Aave
$ws.Range("D43").Value = '95.67'
$ws.Range("E43").Value = '  -1.43%  '

# HuobiToken
$ws.Range("E44").Value = '  -0.40%  '

# Cronos
$ws.Range("D45").Value = '0.0911'
$ws.Range("E45").Value = '  -1.16%  '

# TrustWalletToken
$ws.Range("E46").Value = '  -2.22%  '

# FTXToken
$ws.Range("D47").Value = '4.01'
$ws.Range("E47").Value = '  -1.94%  '

# ARBITRUM
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -1.32%  '

# MXToken
$ws.Range("E49").Value = '  -0.21%  '

# FraxShare
$ws.Range("D50").Value = '7.08'
$ws.Range("E50").Value = '  +0.46%  '

# RocketPoolETH
$ws.Range("D51").Value = '2.215.50'
$ws.Range("E51").Value = '  -1.17%  '

$rng.Style = "Normal"
